# "Adding Header field in to 1PProfile"
#
# The STATUS column (L) test-result data ("PASS" in every data row) is
# removed, along with a handful of stray empty placeholder cells in the
# QUERYSTRING (G), DEPENDENCYTESTS (I) and STORE (K) columns that carried
# no value and no formatting. Cells that still carry formatting (a style)
# or real content are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (QUERYSTRING): drop the bare/unstyled empty cells.
$ws.Range("G2").ClearContents() | Out-Null
$ws.Range("G7:G12").ClearContents() | Out-Null
$ws.Range("G14:G16").ClearContents() | Out-Null

# Column I (DEPENDENCYTESTS): drop the two bare/unstyled empty cells.
$ws.Range("I2").ClearContents() | Out-Null
$ws.Range("I7").ClearContents() | Out-Null

# Column K (STORE): drop the bare/unstyled empty cells.
$ws.Range("K7:K12").ClearContents() | Out-Null
$ws.Range("K14:K16").ClearContents() | Out-Null

# Column L (STATUS): clear all of the "PASS" data, keep only the header.
$ws.Range("L2:L16").ClearContents() | Out-Null

# Reflect the cleared range as the active selection.
$ws.Range("L2:L16").Select() | Out-Null
